# Update the "想去人数" (interested-people count) figures for a handful of
# conventions across the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 11999
$wsExhibition.Range("F7").Value = 227
$wsExhibition.Range("F8").Value = 11900
$wsExhibition.Range("F14").Value = 5893

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 11999
$wsAll.Range("F10").Value = 227
$wsAll.Range("F11").Value = 11900
$wsAll.Range("F18").Value = 5893
